$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 2 = Liz's standup entry: fill in what she did, what she plans to do,
# and her obstacles.
$cell = $t.Cell(2, 2)
$cell.Range.Text = "Created Django project began creating user and profile logic"

$cell = $t.Cell(2, 3)
$cell.Range.Text = "Finish Users and profile logic"

$cell = $t.Cell(2, 4)
$cell.Range.Text = "I" + [char]0x2019 + "m still getting used to Django. "
$rng = $cell.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()

# Row 3 = Ian's standup entry: the "plan to do" column text is unchanged in
# content but is split across two runs in the authored edit.
$cell = $t.Cell(3, 3)
$cell.Range.Text = "Create list of tools and "
$rng = $cell.Range
$rng.Collapse(0)
$rng.InsertAfter("compile images and descriptions for those tools.")
